$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old second data row (row 3) ---
$ws.Rows("3").Delete()

# --- Update header row (row 1) with the new, expanded set of columns ---
$ws.Range("A1").Value = "Hospital"
$ws.Range("B1").Value = "CSN"
$ws.Range("C1").Value = "Patient"
$ws.Range("D1").Value = "SEX_NAME"
$ws.Range("E1").Value = "Age_At_Admission"
$ws.Range("F1").Value = "DAYS_TO_READMISSION"
$ws.Range("G1").Value = "Mortality from Index Culture (Days)"
$ws.Range("H1").Value = "BMI"
$ws.Range("I1").Value = "LOS_Days"
$ws.Range("J1").Value = "Admission Date"
$ws.Range("K1").Value = "Discharge Date"
$ws.Range("L1").Value = "ICU_Encounter"
$ws.Range("M1").Value = "Culture"
$ws.Range("N1").Value = "Index_Culture"
$ws.Range("O1").Value = "Polymicrobial_Infection"
$ws.Range("P1").Value = "Patient_Had_ID_consult"
$ws.Range("Q1").Value = "Diagnosis_Endocarditis_This_Encounter"
$ws.Range("R1").Value = "Diagnosis_Osteomyelitis_This_Encounter"
$ws.Range("S1").Value = "Cefepime"
$ws.Range("T1").Value = "Piperacillin/Tazobactam"
$ws.Range("U1").Value = "Group"
$ws.Range("V1").Value = "Total_DOT"
$ws.Range("W1").Value = "Last_Admin"

# Apply the existing header style (bold font + border + centered alignment) to all the
# newly-introduced header cells so the whole row 1 looks consistent.
$ws.Range("A1").Copy()
$ws.Range("D1:W1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Rebuild the single remaining data row (row 2) under the new column layout ---
# Clear out anything left over from the previous layout (values AND formats, so no stray
# number formatting lingers on cells that should now be plain/blank).
$ws.Range("A2:W2").ClearContents()
$ws.Range("A2:W2").ClearFormats()

$ws.Range("B2").Value = 1
$ws.Range("S2").Value = 8
$ws.Range("T2").Value = 7
$ws.Range("U2").Value = 1
$ws.Range("V2").Value = 12
$ws.Range("W2").Value = 44598

# Give the Last_Admin cell the same date-style formatting used previously (numFmt 165).
$ws.Range("W2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
